# Apply "Append: 2026-02-12 06:59 JST" update to the "ランサーズ" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-02-12 06:59:19"

# --- Update row 2 (timestamp only; all other columns unchanged) ---
$ws.Range("A2").Value = $newTimestamp

# --- Update row 3 ---
$ws.Range("A3").Value = $newTimestamp
$ws.Range("B3").Value = "【急募】ノーコードで実現するLINE×AI恋愛体験サービスMVP開発"
$ws.Range("D3").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5490408"
$ws.Range("G3").Value = 368
# H3 unchanged ("🔥AI,Ai ◆開発")

# --- Update row 4 ---
$ws.Range("A4").Value = $newTimestamp
$ws.Range("B4").Value = "自動化システム"
$ws.Range("D4").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5477084"
$ws.Range("G4").Value = 110
$ws.Range("H4").Value = "◆自動化"

# --- Update row 5 ---
$ws.Range("A5").Value = $newTimestamp
$ws.Range("B5").Value = "【急募】FileMakerシステムのデバッグとレイアウト修正依頼"
$ws.Range("D5").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5490478"
$ws.Range("G5").Value = 28
$ws.Range("H5").ClearContents()

# --- Update row 6 ---
$ws.Range("A6").Value = $newTimestamp
$ws.Range("B6").Value = "【長期】寝具ブランドのAmazon・楽天市場 運用代行パートナー募集"
$ws.Range("D6").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5486471"
$ws.Range("G6").Value = 25
$ws.Range("H6").ClearContents()

# --- Update row 7 ---
$ws.Range("A7").Value = $newTimestamp
$ws.Range("B7").Value = "プロジェクトマネジメント"
$ws.Range("D7").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5490062"
$ws.Range("G7").Value = 25
$ws.Range("H7").ClearContents()

# --- Update row 8 ---
$ws.Range("A8").Value = $newTimestamp
$ws.Range("B8").Value = "【SES経営者向け】事業立ち上げについてお話をお伺いできる方を募集します!"
$ws.Range("D8").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5490407"
$ws.Range("G8").Value = 10
$ws.Range("H8").ClearContents()

# --- Drop the now-stale rows 9-15 ---
$ws.Rows("9:15").Delete()

# --- Rebuild hyperlinks for F2:F8 so the relationship targets match the
#     refreshed URLs (row delete / value assignment do not retarget the
#     existing hyperlink relationships automatically). ---
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5489981")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5490408")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5477084")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5490478")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5486471")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5490062")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5490407")
$ws.Range("F2:F8").Style = "Hyperlink"

# --- Column B narrows from 52 to 39 character-width units.
#     Excel's COM ColumnWidth property is expressed in units that are offset
#     from the raw OOXML <col width> by ~5/6 of a character (the default
#     column padding); subtract that offset so the saved width is exactly 39. ---
$ws.Columns("B").ColumnWidth = 39 - 5/6
